# Add newly-discovered Italian wines to the flash-card workbook.
$wb = $excel.ActiveWorkbook

# --- "European design. & varieties" sheet: two new appellation/grape rows ---
$ws5 = $wb.Worksheets.Item("European design. & varieties")
$ws5.Range("A56").Value = "Lombardie (Red), Valtellina, Italy"
$ws5.Range("B56").Value = "Nebbiolo"
$ws5.Range("A57").Value = "Latium (White)"
$ws5.Range("B57").Value = "Frascati"

# Scroll/select as the author left the sheet after adding the rows.
$ws5.Activate()
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$ws5.Range("A58").Select() | Out-Null

# --- "Wine to discover (todo)" sheet: one new todo entry ---
$ws6 = $wb.Worksheets.Item("Wine to discover (todo)")
$ws6.Range("A2").Value = "White wines"
$ws6.Range("B2").Value = 5357
$ws6.Range("C2").Value = "Kindle"
$ws6.Range("D2").Value = "Region / Appellation"
